$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.534264326095581
$ws.Range("B1").Value = 2.341079473495483
$ws.Range("C1").Value = 4.278539180755615
$ws.Range("D1").Value = 1.837785482406616
$ws.Range("E1").Value = 0.819733202457428
